$wb = $excel.ActiveWorkbook
$sheet3 = $wb.Worksheets.Item("Sheet3")
$sheet3.Name = "tip deflection results"
$sheet3.Range("A1").Value = "test #"
$sheet3.Range("A2").Value = 7
$sheet3.Range("A3").Value = 14
$sheet3.Range("A4").Value = 21
$sheet3.Range("A5").Value = 22
$sheet3.Range("D1").Value = "constload tip deflection [m]"
$sheet3.Range("D2").Value = 0.4682401240911
$sheet3.Range("D3").Value = 0.48980424001869999
$sheet3.Range("D4").Value = 0.54492962852030002
$sheet3.Range("D5").Value = 0.56385624898730002

$chartObj = $sheet3.Shapes.AddChart2(-1, 51)
$chart = $chartObj.Chart
$chart.SeriesCollection.NewSeries()
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(`'tip deflection results`'!`$A`$1,`'tip deflection results`'!`$A`$2:`$A`$5,`'tip deflection results`'!`$D`$2:`$D`$5,1)"

$chart.HasTitle = $false

$catAxis = $chart.Axes(1)
$catAxis.HasTitle = $true
$catAxis.AxisTitle.Text = "test #"

$valAxis = $chart.Axes(2)
$valAxis.HasTitle = $true
$valAxis.AxisTitle.Text = "constload tip deflection [m]"
$valAxis.MinimumScale = 0.3
$valAxis.TickLabels.NumberFormat = "0.00"
$valAxis.HasMajorGridlines = $true

$chart.HasLegend = $false
